$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Updated "Strata weight" (column N) values for the manuscript table/model
$values = @{
    3  = 0.33
    4  = 0.33
    5  = 0.13
    6  = 0.33
    7  = 0.33
    8  = 0.13
    9  = 0.22
    10 = 0.22
    11 = 0.33
    12 = 0.33
    13 = 0.33
    14 = 0.22
    15 = 0.22
    16 = 0.22
    17 = 0.13
    18 = 0.13
    19 = 5.01
    20 = 0.13
    21 = 0.13
    22 = 0.27
    23 = 0.27
    24 = 0.33
    25 = 0.27
    26 = 0.27
    27 = 0.27
    28 = 0.27
    29 = 0.13
    30 = 0.27
    31 = 1.53
    32 = 0.27
    33 = 0.13
    34 = 0.13
    35 = 0.13
    36 = 0.27
}

foreach ($row in $values.Keys) {
    $ws.Range("N$row").Value = $values[$row]
}

# Reflect the last-edited cell selection as in the authored workbook
$ws.Activate()
$ws.Range("N35").Select()
